# Generate Report for Handback
#
# The nightly localization-status report is being re-generated now that the
# handback step has actually run: files that used to be "Ready for handoff"
# are now "Handed back: in sync with en-US", the per-language sheets grow a
# "Latest Target File" / "Latest Handback File" column pair (with real
# hyperlinks, mirroring the existing Source/Handoff-file links), and the
# "Latest Handback DateTime" column picks up a real timestamp instead of the
# 0001-01-01 placeholder.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$sourceRepoBase = "https://github.com/OpenLocalizationTest/oltest/blob/5de96a0e3c72f4af2c7ff4b6a3a0d5a1950f0360/e2e"
$zhHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/820165064aa34e84e813693a9a6ec975ebdef250/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/37e9cff96110580122ba9d1886e15de7b93d12f3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

function Update-StatusColumn($ws, $cells) {
    foreach ($addr in $cells) {
        $rng = $ws.Range($addr)
        if ($rng.Value2 -eq $oldStatus) {
            $rng.Value = $newStatus
        }
    }
}

# --- Overview sheet: Status cells for both languages, both rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
Update-StatusColumn $wsOverview @("B2", "C2", "B3", "C3")

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-StatusColumn $wsZh @("C2", "C3")

# Row 2 (a.md): Latest Target File (F) / Latest Handback File (G)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "$sourceRepoBase/a.md", "", "", "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhHandoffUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

# Row 3 (b.md): same target/handback file pair as row 2
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "$sourceRepoBase/a.md", "", "", "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhHandoffUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

# Latest Handback DateTime now has a real timestamp
$wsZh.Range("H2").Value = "2016-03-18 10:24:26"
$wsZh.Range("H3").Value = "2016-03-18 10:24:26"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
Update-StatusColumn $wsDe @("C2", "C3")

# Row 2 (a.md): Latest Target File (F) / Latest Handback File (G)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "$sourceRepoBase/a.md", "", "", "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deHandoffUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

# Row 3 (b.md): same target/handback file pair as row 2
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "$sourceRepoBase/a.md", "", "", "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deHandoffUrl, "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

# Latest Handback DateTime now has a real timestamp (de-de ran a little later)
$wsDe.Range("H2").Value = "2016-03-18 10:24:31"
$wsDe.Range("H3").Value = "2016-03-18 10:24:31"
